# Scoreboard.xlsx update — "Add files via upload"
#
# Adds the Minute3 / Second3 (and, on ScoreM, Rep3-derived "points") columns
# N/O/P for each team row on the ScoreM and ScoreF sheets, fixes a data
# entry (J23 on ScoreF), and updates which sheet/cell is active & selected.

$wb = $excel.ActiveWorkbook

$wsM = $wb.Worksheets.Item("ScoreM")
$wsF = $wb.Worksheets.Item("ScoreF")

# ---------------------------------------------------------------------
# ScoreM: columns N (Minute3), O (Second3), P (points) for rows 2-24
# ---------------------------------------------------------------------
$scoreMRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24)
$scoreMN    = @(4,8,8,8,8,4,8,8,8,8,5,8,8,5,5,8,8,6,8,5,8,6,5)
$scoreMO    = @(10,0,0,0,0,48,0,0,0,0,41,0,0,44,5,0,0,29,0,3,0,18,37)
$scoreMP    = @(156,148,151,136,142,156,122,124,128,152,156,136,136,156,156,6,150,136,138,156,143,156,156)

for ($i = 0; $i -lt $scoreMRows.Length; $i++) {
    $r = $scoreMRows[$i]
    $wsM.Cells.Item($r, 14).Value = $scoreMN[$i]
    $wsM.Cells.Item($r, 15).Value = $scoreMO[$i]
    $wsM.Cells.Item($r, 16).Value = $scoreMP[$i]
}

# ---------------------------------------------------------------------
# ScoreF: columns N (Minute3), O (Second3) for rows 2-25 (no points yet)
# ---------------------------------------------------------------------
for ($r = 2; $r -le 25; $r++) {
    $wsF.Cells.Item($r, 14).Value = 8
    $wsF.Cells.Item($r, 15).Value = 0
}

# Data fix: Rep2 for row 23 on ScoreF (76 -> 85); the M23 "points" formula
# (SUM(I23:L23)) recalculates automatically.
$wsF.Range("J23").Value = 85

# ---------------------------------------------------------------------
# Active sheet / selection changes
#  - ScoreF loses tabSelected, its selection moves to J24
#  - ScoreM becomes the active/selected tab, with selection at P5
# ---------------------------------------------------------------------
$wsF.Range("J24").Select()
$wsM.Activate()
$wsM.Range("P5").Select()
